$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing Text storage (matches source
# workbook, which stores every cell - including numeric-looking prices -
# as inline text) and then restoring the default "Normal" style so no
# stray number-format style gets attached to the cell.
function Set-TextValue {
    param($Sheet, [string]$Ref, [string]$Val)
    $Sheet.Range($Ref).NumberFormat = "@"
    $Sheet.Range($Ref).Value = $Val
    $Sheet.Range($Ref).Style = "Normal"
}

Set-TextValue $ws "D2" "69.375.84"
$ws.Range("E2").Value = "  +1.75%  "
Set-TextValue $ws "D3" "3.908.07"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws "D5" "529.26"
$ws.Range("E5").Value = "  +9.48%  "
Set-TextValue $ws "D6" "143.99"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue $ws "D9" "0.718"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("E11").Value = "  -4.82%  "
Set-TextValue $ws "D12" "42.01"
$ws.Range("E12").Value = "  -2.80%  "
Set-TextValue $ws "D13" "4.538.09"
$ws.Range("E13").Value = "  +0.76%  "
Set-TextValue $ws "D14" "10.25"
$ws.Range("E14").Value = "  -2.40%  "
Set-TextValue $ws "D15" "3.931.67"
$ws.Range("E15").Value = "  +0.57%  "
Set-TextValue $ws "D16" "1.24"
$ws.Range("E16").Value = "  +9.86%  "
Set-TextValue $ws "D17" "14.00"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("E18").Value = "  -0.74%  "
Set-TextValue $ws "D19" "19.74"
$ws.Range("E19").Value = "  -1.01%  "
Set-TextValue $ws "D20" "69.281.10"
$ws.Range("E20").Value = "  +1.58%  "
Set-TextValue $ws "D21" "427.19"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  -5.52%  "
Set-TextValue $ws "D23" "88.59"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  -4.58%  "
Set-TextValue $ws "D25" "4.06"
$ws.Range("E25").Value = "  +10.92%  "
Set-TextValue $ws "D26" "11.49"
$ws.Range("E26").Value = "  -6.37%  "
Set-TextValue $ws "D27" "10.61"
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("E28").Value = "  -2.26%  "
Set-TextValue $ws "D29" "13.13"
$ws.Range("E29").Value = "  -2.71%  "
Set-TextValue $ws "D30" "674.40"
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("E32").Value = "  -2.91%  "
Set-TextValue $ws "D33" "68.79"
$ws.Range("E33").Value = "  +11.98%  "
Set-TextValue $ws "D34" "0.0₃0888"
$ws.Range("E34").Value = "  +0.66%  "
Set-TextValue $ws "D35" "0.436"
$ws.Range("E35").Value = "  +10.55%  "
Set-TextValue $ws "D36" "5.96"
$ws.Range("E36").Value = "  -1.49%  "
Set-TextValue $ws "D37" "40.03"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("E38").Value = "  +1.87%  "
Set-TextValue $ws "D39" "0.998"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.01%  "
Set-TextValue $ws "D41" "3.21"
$ws.Range("E41").Value = "  +4.67%  "
Set-TextValue $ws "D42" "0.0481"
$ws.Range("E42").Value = "  -3.66%  "
Set-TextValue $ws "D43" "3.16"
$ws.Range("E43").Value = "  +7.16%  "
Set-TextValue $ws "D44" "2.81"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws "D45" "0.000294"
$ws.Range("E45").Value = "  +22.46%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D46" "3.40"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D47" "0.141"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D48" "2.99"
$ws.Range("E48").Value = "  +6.78%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D49" "0.0₆0351"
$ws.Range("E49").Value = "  +5.57%  "
Set-TextValue $ws "D50" "2.744.87"
$ws.Range("E50").Value = "  +13.70%  "
Set-TextValue $ws "D51" "144.82"
$ws.Range("E51").Value = "  +0.29%  "
